# Update team-specific transition-matrix probabilities on Sheet1
# with refreshed values from games pulled March 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1853146853146853
$ws.Range("C2").Value = 0.5734265734265734
$ws.Range("J2").Value = 0.01048951048951049
$ws.Range("P2").Value = 0.1328671328671329
$ws.Range("S2").Value = 0.0979020979020979
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.06741573033707865
$ws.Range("J3").Value = 0.03932584269662921
$ws.Range("P3").Value = 0.7808988764044944
$ws.Range("S3").Value = 0.1067415730337079
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.1012658227848101
$ws.Range("D6").Value = 0.008438818565400843
$ws.Range("F6").Value = 0.08860759493670886
$ws.Range("J6").Value = 0.1814345991561181
$ws.Range("O6").Value = 0.02531645569620253
$ws.Range("Q6").Value = 0.1518987341772152
$ws.Range("R6").Value = 0.06751054852320675
$ws.Range("S6").Value = 0.3755274261603376
$ws.Range("B7").Value = 0.08839779005524862
$ws.Range("D7").Value = 0.03314917127071823
$ws.Range("F7").Value = 0.03314917127071823
$ws.Range("J7").Value = 0.09944751381215469
$ws.Range("O7").Value = 0.01104972375690608
$ws.Range("Q7").Value = 0.1933701657458564
$ws.Range("R7").Value = 0.0718232044198895
$ws.Range("S7").Value = 0.4696132596685083
$ws.Range("B8").Value = 0.1047957371225577
$ws.Range("D8").Value = 0.02131438721136767
$ws.Range("E8").Value = 0.003552397868561279
$ws.Range("F8").Value = 0.07282415630550622
$ws.Range("J8").Value = 0.1119005328596803
$ws.Range("O8").Value = 0.01243339253996448
$ws.Range("Q8").Value = 0.1705150976909414
$ws.Range("R8").Value = 0.07460035523978685
$ws.Range("S8").Value = 0.4280639431616341
$ws.Range("B9").Value = 0.09389671361502347
$ws.Range("D9").Value = 0.02347417840375587
$ws.Range("F9").Value = 0.06572769953051644
$ws.Range("J9").Value = 0.1126760563380282
$ws.Range("O9").Value = 0.02347417840375587
$ws.Range("Q9").Value = 0.1596244131455399
$ws.Range("R9").Value = 0.07042253521126761
$ws.Range("S9").Value = 0.4507042253521127
$ws.Range("B10").Value = 0.09154383242823895
$ws.Range("D10").Value = 0.01318851823118697
$ws.Range("F10").Value = 0.06904577191621412
$ws.Range("J10").Value = 0.1256788207913111
$ws.Range("O10").Value = 0.01939487975174554
$ws.Range("Q10").Value = 0.2063615205585725
$ws.Range("R10").Value = 0.07680372381691233
$ws.Range("S10").Value = 0.3979829325058185
$ws.Range("G11").Value = 0.1578947368421053
$ws.Range("J11").Value = 0.07236842105263158
$ws.Range("K11").Value = 0.2138157894736842
$ws.Range("L11").Value = 0.5230263157894737
$ws.Range("S11").Value = 0.03289473684210526
$ws.Range("G12").Value = 0.6761363636363636
$ws.Range("J12").Value = 0.1761363636363636
$ws.Range("K12").Value = 0.02272727272727273
$ws.Range("L12").Value = 0.04545454545454546
$ws.Range("S12").Value = 0.07954545454545454
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.15
$ws.Range("F15").Value = 0.0193050193050193
$ws.Range("H15").Value = 0.1891891891891892
$ws.Range("I15").Value = 0.07335907335907337
$ws.Range("J15").Value = 0.3590733590733591
$ws.Range("K15").Value = 0.0694980694980695
$ws.Range("M15").Value = 0.007722007722007722
$ws.Range("O15").Value = 0.09266409266409266
$ws.Range("S15").Value = 0.1891891891891892
$ws.Range("F16").Value = 0.005128205128205128
$ws.Range("H16").Value = 0.2564102564102564
$ws.Range("I16").Value = 0.08205128205128205
$ws.Range("J16").Value = 0.3641025641025641
$ws.Range("K16").Value = 0.09743589743589744
$ws.Range("M16").Value = 0.01538461538461539
$ws.Range("N16").Value = 0.005128205128205128
$ws.Range("O16").Value = 0.06153846153846154
$ws.Range("S16").Value = 0.1128205128205128
$ws.Range("F17").Value = 0.01515151515151515
$ws.Range("H17").Value = 0.1861471861471861
$ws.Range("I17").Value = 0.1103896103896104
$ws.Range("J17").Value = 0.3852813852813853
$ws.Range("K17").Value = 0.07575757575757576
$ws.Range("M17").Value = 0.01298701298701299
$ws.Range("O17").Value = 0.05627705627705628
$ws.Range("S17").Value = 0.158008658008658
$ws.Range("F18").Value = 0.01081081081081081
$ws.Range("H18").Value = 0.1621621621621622
$ws.Range("I18").Value = 0.0918918918918919
$ws.Range("J18").Value = 0.4216216216216216
$ws.Range("K18").Value = 0.07027027027027027
$ws.Range("M18").Value = 0.01621621621621622
$ws.Range("O18").Value = 0.05945945945945946
$ws.Range("S18").Value = 0.1675675675675676
$ws.Range("F19").Value = 0.01739130434782609
$ws.Range("H19").Value = 0.2307692307692308
$ws.Range("I19").Value = 0.06822742474916388
$ws.Range("J19").Value = 0.345819397993311
$ws.Range("K19").Value = 0.09565217391304348
$ws.Range("M19").Value = 0.01672240802675585
$ws.Range("O19").Value = 0.07290969899665552
$ws.Range("S19").Value = 0.1525083612040134
